$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entries for 周四 (2019年5月16日) - three new rows in columns C/D
$ws.Range("C53").Value = "崔梦婷Dao单表建立成功，TEST单表通过"
$ws.Range("D53").Value = "19:30--19:40"

$ws.Range("C54").Value = "邢朋举Service层通过，缺少TEST"
$ws.Range("D54").Value = "20:00--20:20"

$ws.Range("D55").Value = "20:30--21:00"
$ws.Range("C55").Value = "李博文DaoTest AND ServiceTest均通过"

$ws.Range("D55").Select()
